# Rebuild the recipients data sheet with the new columns (G:J) and the
# updated/re-ordered rows (including two brand-new rows coming from the
# new web-based import/export flow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Due Amount"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Receipt No"
$ws.Range("G1").Value = "Recipient Email"
$ws.Range("H1").Value = "Payment_Status"
$ws.Range("I1").Value = "Amount_Paid"
$ws.Range("J1").Value = "Last Updated"

# ---- Row 2 : SARWAR ----
$ws.Range("A2").Value = "SARWAR"
$ws.Range("B2").Value = 2500
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = "'2026-02-02"
$ws.Range("E2").Value = "Dummy payment"
$ws.Range("F2").Value = 2001
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# ---- Row 3 : Jane Smith ----
$ws.Range("A3").Value = "Jane Smith"
$ws.Range("B3").Value = 3500
$ws.Range("C3").Value = 1500
$ws.Range("D3").Value = "'2026-01-30"
$ws.Range("E3").Value = "Partial payment"
$ws.Range("F3").Value = 2003
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""

# ---- Row 4 : Test Payment (new row from the web app test) ----
$ws.Range("A4").Value = "Test Payment"
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = "2026-02-01T15:00:00.000Z"
$ws.Range("E4").Value = "Test payment for verification"
$ws.Range("F4").Value = 9999
$ws.Range("G4").Value = "sarwarofficial2006@gmail.com"
$ws.Range("H4").Value = "DUE"
$ws.Range("I4").Value = "¥1500"
$ws.Range("J4").Value = "2026-02-01T17:34:33.169Z"

# ---- Row 5 : Maria Garcia ----
$ws.Range("A5").Value = "Maria Garcia"
$ws.Range("B5").Value = 1800
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "'2026-01-28"
$ws.Range("E5").Value = "Retainer"
$ws.Range("F5").Value = 2005
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""

# ---- Row 6 : Ali (new row) ----
$ws.Range("A6").Value = "Ali"
$ws.Range("B6").Value = 2500
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = "2026-01-31T15:00:00.000Z"
$ws.Range("E6").Value = "Dummy payment"
$ws.Range("F6").Value = 2006
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""

# ---- Row 7 : John Doe ----
$ws.Range("A7").Value = "John Doe"
$ws.Range("B7").Value = 5000
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "'2026-02-01"
$ws.Range("E7").Value = "Service fee"
$ws.Range("F7").Value = 2002
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""

# ---- Row 8 : Ahmed Khan ----
$ws.Range("A8").Value = "Ahmed Khan"
$ws.Range("B8").Value = 4200
$ws.Range("C8").Value = 4200
$ws.Range("D8").Value = "'2026-02-02"
$ws.Range("E8").Value = "Invoice #001"
$ws.Range("F8").Value = 2004
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
